# Apply "changes to rest api document" to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# C5: clarify merge condition to also include "property".
$ws.Range("C5").Value = "merge if data is entity or property"

# Move the sheet's active selection from B8 down to A9 (one row below the
# table, column A) - matches the new <selection activeCell="A9" sqref="A9"/>.
$ws.Range("A9").Select() | Out-Null

# Best-effort: reflect the new saved window geometry. (Harmless if the
# headless host doesn't expose/persist this cosmetic window-chrome size.)
try {
    $win = $wb.Windows.Item(1)
    $win.Width = 18660
    $win.Height = 8080
} catch {
}
